# Update "想去人数" (want-to-go count) figures in the 展览 sheet and the
# corresponding rows mirrored into the 全部类型 sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) -------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 367   # was 366
$wsExpo.Range("F3").Value = 786   # was 783
$wsExpo.Range("F4").Value = 271   # was 272
$wsExpo.Range("F5").Value = 850   # was 843
$wsExpo.Range("F6").Value = 2098  # was 2083
$wsExpo.Range("F7").Value = 188   # was 187

# --- Sheet "全部类型" (all types, mirrors the other sheets' rows) --------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 367    # was 366
$wsAll.Range("F3").Value = 786    # was 783
$wsAll.Range("F4").Value = 271    # was 272
$wsAll.Range("F7").Value = 850    # was 843
$wsAll.Range("F8").Value = 2098   # was 2083
$wsAll.Range("F10").Value = 188   # was 187

$wb.Save()
